# "Fruta / hortaliza, semanal"
#
# The published dataset gained one new weekly price observation for
# Rabanito (Vega Central Mapocho de Santiago). It sits at the top of the
# date-descending price history, so a new row is inserted at row 59
# (pushing the existing rows 59-186 down to 60-187) and populated with
# the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59; this shifts the existing rows 59..186
# down to 60..187, carrying their values/formatting with them.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new observation.
$ws.Cells.Item(59, 1).Value = 9
$ws.Cells.Item(59, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(59, 3).Value = "Metropolitana"
$ws.Cells.Item(59, 4).Value = 44519
$ws.Cells.Item(59, 5).Value = 13
$ws.Cells.Item(59, 6).Value = 300000001
$ws.Cells.Item(59, 7).Value = "Rabanito"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 7900
$ws.Cells.Item(59, 11).Value = 2500
$ws.Cells.Item(59, 12).Value = 3000
$ws.Cells.Item(59, 13).Value = 2747
$ws.Cells.Item(59, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(59, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(59, 16).Value = 27
$ws.Cells.Item(59, 17).Value = 100
$ws.Cells.Item(59, 18).Value = "Hortaliza"
